$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1875
$ws.Range("J17").Value = 1875
$ws.Range("L17").Value = 5625
$ws.Range("N17").Value = -5961
$ws.Range("H96").Value = 1781.6428
$ws.Range("I96").Value = 782.1667
$ws.Range("K96").Value = 2346.5001
$ws.Range("M96").Value = -973.5001000000002
$ws.Range("H112").Value = 1687.7
$ws.Range("J112").Value = 1836
$ws.Range("L112").Value = 5508
$ws.Range("N112").Value = -7724
$ws.Range("H116").Value = 5747.4707
$ws.Range("I116").Value = 5518.273
$ws.Range("J116").Value = 6167.6665
$ws.Range("K116").Value = 5518.273
$ws.Range("L116").Value = 6167.6665
$ws.Range("M116").Value = -2076.273
$ws.Range("N116").Value = -13051.6665
$ws.Range("H125").Value = 2718.25
$ws.Range("I125").Value = 2458.1667
$ws.Range("J125").Value = 3498.5
$ws.Range("K125").Value = 22123.5003
$ws.Range("L125").Value = 31486.5
$ws.Range("M125").Value = -19663.5003
$ws.Range("N125").Value = -36406.5
$ws.Range("H132").Value = 2844.7856
$ws.Range("I132").Value = 1925.2222
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 5775.6666
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3245.6666
$ws.Range("N132").Value = -18560

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 849.2727
$ws.Range("J2").Value = 1259.6
$ws.Range("L2").Value = 1259.6
$ws.Range("N2").Value = -1485.6
$ws.Range("H45").Value = 13159944
$ws.Range("I45").Value = 17858746
$ws.Range("J45").Value = 3297.8
$ws.Range("K45").Value = 17858746
$ws.Range("L45").Value = 3297.8
$ws.Range("M45").Value = -17858369
$ws.Range("N45").Value = -4051.8
$ws.Range("H97").Value = 946.6
$ws.Range("I97").Value = 855.9048
$ws.Range("J97").Value = 1422.75
$ws.Range("K97").Value = 855.9048
$ws.Range("L97").Value = 1422.75
$ws.Range("M97").Value = -359.9048
$ws.Range("N97").Value = -2414.75
$ws.Range("H102").Value = 11251.615
$ws.Range("I102").Value = 8976.909
$ws.Range("J102").Value = 23762.5
$ws.Range("K102").Value = 8976.909
$ws.Range("L102").Value = 23762.5
$ws.Range("M102").Value = -7354.909
$ws.Range("N102").Value = -27006.5
$ws.Range("H110").Value = 4050.5833
$ws.Range("I110").Value = 3972.875
$ws.Range("K110").Value = 3972.875
$ws.Range("M110").Value = -1927.875
$ws.Range("H116").Value = 849.2727
$ws.Range("J116").Value = 1259.6
$ws.Range("L116").Value = 1259.6
$ws.Range("N116").Value = -5847.6
$ws.Range("H122").Value = 4177.727
$ws.Range("I122").Value = 1987
$ws.Range("J122").Value = 4999.25
$ws.Range("K122").Value = 5961
$ws.Range("L122").Value = 14997.75
$ws.Range("M122").Value = -3511
$ws.Range("N122").Value = -19897.75
$ws.Range("H127").Value = 39999
$ws.Range("J127").Value = 39999
$ws.Range("L127").Value = 39999
$ws.Range("N127").Value = -49919
$ws.Range("H132").Value = 2761.8484
$ws.Range("I132").Value = 2785.4375
$ws.Range("K132").Value = 8356.3125
$ws.Range("M132").Value = -5826.3125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 849.2727
$ws.Range("J3").Value = 1259.6
$ws.Range("L3").Value = 1259.6
$ws.Range("N3").Value = -1487.6
$ws.Range("H94").Value = 729.35297
$ws.Range("I94").Value = 568.1667
$ws.Range("J94").Value = 1116.2
$ws.Range("K94").Value = 568.1667
$ws.Range("L94").Value = 1116.2
$ws.Range("M94").Value = -117.1667
$ws.Range("N94").Value = -2018.2
$ws.Range("H99").Value = 2937.4736
$ws.Range("I99").Value = 2088.9
$ws.Range("J99").Value = 3880.3333
$ws.Range("K99").Value = 2088.9
$ws.Range("L99").Value = 3880.3333
$ws.Range("M99").Value = -590.9000000000001
$ws.Range("N99").Value = -6876.3333
$ws.Range("H105").Value = 2928.8
$ws.Range("I105").Value = 2541.1428
$ws.Range("K105").Value = 2541.1428
$ws.Range("M105").Value = -794.1428000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2005199.2
$ws.Range("I62").Value = 3337998.8
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 3337998.8
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3337374.8
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 2005199.2
$ws.Range("I65").Value = 3337998.8
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 16689994
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -16686874
$ws.Range("N65").Value = -36240
$ws.Range("H105").Value = 1461.5652
$ws.Range("I105").Value = 1403.1111
$ws.Range("J105").Value = 1672
$ws.Range("K105").Value = 1403.1111
$ws.Range("L105").Value = 1672
$ws.Range("M105").Value = 343.8888999999999
$ws.Range("N105").Value = -5166
$ws.Range("H107").Value = 1891.6666
$ws.Range("I107").Value = 831.6429000000001
$ws.Range("J107").Value = 3033.2307
$ws.Range("K107").Value = 831.6429000000001
$ws.Range("L107").Value = 3033.2307
$ws.Range("M107").Value = 1088.3571
$ws.Range("N107").Value = -6873.2307

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2275.5833
$ws.Range("I12").Value = 2973.375
$ws.Range("K12").Value = 8920.125
$ws.Range("M12").Value = -8747.125
$ws.Range("H107").Value = 1261.5714
$ws.Range("J107").Value = 1261.5714
$ws.Range("L107").Value = 3784.7142
$ws.Range("N107").Value = -7624.7142

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10406.471
$ws.Range("I70").Value = 13493.7
$ws.Range("J70").Value = 5996.143
$ws.Range("K70").Value = 13493.7
$ws.Range("L70").Value = 5996.143
$ws.Range("M70").Value = -13223.7
$ws.Range("N70").Value = -6536.143
$ws.Range("H73").Value = 10406.471
$ws.Range("I73").Value = 13493.7
$ws.Range("J73").Value = 5996.143
$ws.Range("K73").Value = 13493.7
$ws.Range("L73").Value = 5996.143
$ws.Range("M73").Value = -12557.7
$ws.Range("N73").Value = -7868.143
$ws.Range("H80").Value = 2881
$ws.Range("I80").Value = 2305.0625
$ws.Range("J80").Value = 5184.75
$ws.Range("K80").Value = 2305.0625
$ws.Range("L80").Value = 5184.75
$ws.Range("M80").Value = -1307.0625
$ws.Range("N80").Value = -7180.75
$ws.Range("H83").Value = 2881
$ws.Range("I83").Value = 2305.0625
$ws.Range("J83").Value = 5184.75
$ws.Range("K83").Value = 11525.3125
$ws.Range("L83").Value = 25923.75
$ws.Range("M83").Value = -6533.3125
$ws.Range("N83").Value = -35907.75
$ws.Range("H97").Value = 1298.7222
$ws.Range("I97").Value = 1441.125
$ws.Range("J97").Value = 159.5
$ws.Range("K97").Value = 1441.125
$ws.Range("L97").Value = 159.5
$ws.Range("M97").Value = -945.125
$ws.Range("N97").Value = -1151.5
$ws.Range("H102").Value = 3766
$ws.Range("I102").Value = 3087.3333
$ws.Range("J102").Value = 4529.5
$ws.Range("K102").Value = 3087.3333
$ws.Range("L102").Value = 4529.5
$ws.Range("M102").Value = -1465.3333
$ws.Range("N102").Value = -7773.5
$ws.Range("H107").Value = 1014.75
$ws.Range("I107").Value = 988.5
$ws.Range("J107").Value = 1041
$ws.Range("K107").Value = 988.5
$ws.Range("L107").Value = 1041
$ws.Range("M107").Value = 931.5
$ws.Range("N107").Value = -4881
$ws.Range("H132").Value = 34490268
$ws.Range("I132").Value = 47622428
$ws.Range("K132").Value = 142867284
$ws.Range("M132").Value = -142864754

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 43377.69
$ws.Range("I7").Value = 3214.077
$ws.Range("J7").Value = 83541.30499999999
$ws.Range("K7").Value = 3214.077
$ws.Range("L7").Value = 83541.30499999999
$ws.Range("M7").Value = -3102.077
$ws.Range("N7").Value = -83765.30499999999
$ws.Range("H40").Value = 4760.3335
$ws.Range("I40").Value = 3911.9285
$ws.Range("K40").Value = 3911.9285
$ws.Range("M40").Value = -3775.9285
$ws.Range("H120").Value = 145000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 145000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 145000
$ws.Range("M120").Value = $null
$ws.Range("N120").Value = -154676
$ws.Range("H126").Value = 43377.69
$ws.Range("I126").Value = 3214.077
$ws.Range("J126").Value = 83541.30499999999
$ws.Range("K126").Value = 9642.231
$ws.Range("L126").Value = 250623.915
$ws.Range("M126").Value = -7172.231
$ws.Range("N126").Value = -255563.915
$ws.Range("H136").Value = 63610.305
$ws.Range("I136").Value = 9168.866
$ws.Range("J136").Value = 165688
$ws.Range("K136").Value = 27506.598
$ws.Range("L136").Value = 497064
$ws.Range("M136").Value = -24956.598
$ws.Range("N136").Value = -502164

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 50000
$ws.Range("J41").Value = 50000
$ws.Range("L41").Value = 50000
$ws.Range("N41").Value = -50780
$ws.Range("H107").Value = 752.1429000000001
$ws.Range("I107").Value = 751.619
$ws.Range("K107").Value = 2254.857
$ws.Range("M107").Value = -334.857
$ws.Range("H122").Value = 2252.606
$ws.Range("I122").Value = 2322.6296
$ws.Range("K122").Value = 6967.888800000001
$ws.Range("M122").Value = -4517.888800000001
$ws.Range("H132").Value = 180824.23
$ws.Range("I132").Value = 2217.4082
$ws.Range("J132").Value = 1431072
$ws.Range("K132").Value = 6652.2246
$ws.Range("L132").Value = 4293216
$ws.Range("M132").Value = -4122.2246
$ws.Range("N132").Value = -4298276
